# Update countries & provincias Spain
#
# The "Pais" sheet lists countries ranked by total cases (column B, descending).
# In this data refresh, Gabon's case count jumped to 95 - tying it with Aruba
# and moving it up the ranking from its old spot (between Togo and Somalia,
# row 144) to a new spot right after Madagascar (row 135) and before Aruba
# (row 136).
#
# Net effect on the table: a "Gabon" row is inserted after Madagascar (pushing
# Aruba..Togo down by one row), and the old Gabon row further down the list is
# removed (pulling Somalia and everything below back up by one row). Every
# other row/country keeps its own data; only its row position shifts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank row right after "Madagascar" (row 135), before "Aruba"
#    (row 136). This pushes Aruba..Togo (old rows 136-144) down to 137-145.
$ws.Rows(136).Insert()

# 2) The stale "Gabon" entry (old data: 80/0/4/75/0/0/1) has now shifted from
#    row 144 down to row 145. Remove it so Somalia and the rows below move
#    back up by one, restoring the original row count.
$ws.Rows(145).Delete()

# 3) Populate the newly inserted row 136 with Gabon's updated figures.
$ws.Range("A136").Value = "Gabon"
$ws.Range("B136").Value = 95
$ws.Range("C136").Value = 15
$ws.Range("D136").Value = 6
$ws.Range("E136").Value = 88
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 1
